$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report row was inserted as row 14 (Fecha 2021-10-08 /
# serial 44483), pushing the previously-existing rows 14-33 down to 15-34.
$ws.Rows(14).Insert()

$ws.Cells(14,1).Value2 = 10
$ws.Cells(14,2).Value2 = "Vega Modelo de Temuco"
$ws.Cells(14,3).Value2 = "La Araucanía"
$ws.Cells(14,4).Value2 = 44483
$ws.Cells(14,5).Value2 = 9
$ws.Cells(14,6).Value2 = 100112026
$ws.Cells(14,7).Value2 = "Haba"
$ws.Cells(14,8).Value2 = "Sin especificar"
$ws.Cells(14,9).Value2 = "Primera"
$ws.Cells(14,10).Value2 = 80
$ws.Cells(14,11).Value2 = 9000
$ws.Cells(14,12).Value2 = 9000
$ws.Cells(14,13).Value2 = 9000
$ws.Cells(14,14).Value2 = '$/saco 25 kilos'
$ws.Cells(14,15).Value2 = "Provincia de Limarí"
$ws.Cells(14,16).Value2 = 360
$ws.Cells(14,17).Value2 = 25
$ws.Cells(14,18).Value2 = "Hortaliza"
